# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.983.98'
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = '3.102.08'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'526.50"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").Value = "'141.55"
$ws.Range("E6").Value = '  -2.23%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.101.60'
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").Value = "'0.444"
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = "'7.17"
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("E14").Value = '  +3.16%  '
$ws.Range("D15").Value = "'25.65"
$ws.Range("E15").Value = '  -5.78%  '
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("D17").Value = '58.017.90'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '3.100.79'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = "'6.11"
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("D20").Value = "'12.71"
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("D21").Value = "'7.99"
$ws.Range("E21").Value = '  -2.83%  '
$ws.Range("D22").Value = "'343.07"
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D25").Value = "'67.57"
$ws.Range("E25").Value = '  +2.29%  '
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '0.0₃0925'
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = "'6.36"
$ws.Range("E30").Value = '  -7.73%  '
$ws.Range("D31").Value = "'7.29"
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").Value = "'20.99"
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("D34").Value = "'1.19"
$ws.Range("E34").Value = '  -2.93%  '
$ws.Range("D35").Value = "'158.94"
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("D37").Value = "'6.16"
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("D38").Value = "'26.12"
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("E39").Value = '  -4.91%  '
$ws.Range("D40").Value = "'0.0671"
$ws.Range("D41").Value = "'1.59"
$ws.Range("E41").Value = '  +7.70%  '
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("D44").Value = '3.139.90'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("E47").Value = '  +1.60%  '
$ws.Range("D48").Value = '2.268.91'
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").Value = "'20.50"
$ws.Range("E51").Value = '  -2.24%  '
